{"js": "// Remove the \"Heading2\" style from the section-title paragraphs (they become\n// plain body paragraphs again) and refresh the in-text citation markers in\n// the five content paragraphs that discuss the sources, per the commit's\n// citation-check pass.\n\nconst headingTitles = [\"Introduction\", \"Overview of Microfinance\", \"Impact of Microfinance Institutions\", \"Effectiveness in Sustainable Development\", \"Conclusion\"];\n\nconst replacements = [\n  { marker: \"Scholarly work on microfinance reveals a complex l\", after: \"Scholarly work on microfinance reveals a complex landscape of both accomplishments and challenges. Success stories, such as that of the Grameen Bank, underscore the potential of microfinance institutions to empower marginalized populations by improving access to essential services and fostering economic activity within local communities (Smith). However, critiques highlight that while microfinance institutions have achieved economic success in operating microcredit, they often fall short of substantially raising the living standards of their clients (Smith). Furthermore, the literature points to persistent issues such as social divides and gender-based inequalities, which can hinder the full integration of women into economic activities (Smith). These critiques emphasize the need for addressing such challenges to ensure microfinance programs truly benefit the communities they aim to serve, thus promoting both financial and social inclusion (Smith).\" },\n  { marker: \"Microfinance institutions, exemplified by the Gram\", after: \"Microfinance institutions, exemplified by the Grameen Bank, have played a transformative role in empowering local communities through financial inclusion and economic development. The Grameen Bank, renowned for its pioneering approach to microcredit, has significantly improved access to credit for impoverished individuals, thereby enabling them to invest in small-scale businesses and enhance their livelihoods (Ref-u503635). Another notable institution, such as BRAC, complements this impact by providing integrated services that include education and healthcare, further contributing to the well-being of community members (Ref-u503635). These institutions not only elevate household incomes but also stimulate broader economic activity by fostering entrepreneurship and creating employment opportunities in underserved regions (Ref-u503635). As a result, microfinance institutions have become pivotal in driving sustainable development and reducing poverty within local communities, despite facing persistent challenges that require ongoing attention (Ref-u503635).\" },\n  { marker: \"Furthermore, microfinance programs actively allevi\", after: \"Furthermore, microfinance programs actively alleviate poverty through mechanisms like small loans and savings groups, which provide financial access to low-income individuals. These small loans empower recipients to start or expand businesses, leading to increased income and improved living standards (Ref-u026387). Savings groups, on the other hand, enable individuals to pool resources, creating a financial safety net that not only supports personal economic activities but also fosters a culture of saving and financial planning. The collective impact of these mechanisms enhances economic resilience among participants, allowing them to better withstand financial shocks and invest in opportunities that promote long-term stability (Ref-u026387). By integrating these approaches, microfinance programs play a critical role in poverty reduction, thereby contributing to broader economic development goals.\" },\n  { marker: \"However, despite their benefits, microfinance prog\", after: \"However, despite their benefits, microfinance programs face several significant drawbacks that could impede their effectiveness. One major concern is the high interest rates charged by microfinance institutions, which are often necessary to cover operational costs but can burden borrowers and undermine the social mission of these programs (Ref-f532861). Additionally, the dependency created by microfinance services can lead to a cycle where borrowers become reliant on loans without achieving substantial improvements in their living standards (Ref-f532861). Gender-based challenges also persist, as women, who are often the primary beneficiaries of microfinance, may encounter social and cultural barriers that prevent them from fully participating in economic activities (Ref-f532861). These issues highlight the need for a more balanced approach that considers both the financial sustainability of microfinance institutions and the socio-economic empowerment of their clients, ensuring that the programs genuinely contribute to poverty alleviation and economic development.\" },\n  { marker: \"Microfinance has shown considerable potential in p\", after: \"Microfinance has shown considerable potential in promoting sustainable development by improving the economic resilience of disadvantaged populations. Notably, microfinance institutions effectively alleviate credit constraints for impoverished communities, thereby enhancing financial inclusion and facilitating economic empowerment (Ref-u821407). These institutions provide critical financial resources that enable small businesses to flourish, which not only raises household incomes but also stimulates broader economic growth in underdeveloped regions (Ref-u821407). However, despite these positive outcomes, criticisms persist, particularly concerning the high interest rates and potential dependency issues that may limit the long-term impact of microfinance programs (Ref-u821407). Addressing these challenges requires the development of robust regulatory frameworks and policies that support the sustainability of microfinance initiatives while ensuring they contribute meaningfully to poverty reduction and sustainable development goals (Ref-u821407).\" }\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const text = p.text.trim();\n  if (headingTitles.includes(text) && p.style === \"Heading 2\") {\n    p.style = \"Normal\";\n  }\n}\n\nfor (const p of paragraphs.items) {\n  for (const { marker, after } of replacements) {\n    if (p.text.startsWith(marker)) {\n      p.insertText(after, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Heading2\" style from the section-title paragraphs (they become\n# plain body paragraphs again) and refresh the in-text citation markers in\n# the five content paragraphs that discuss the sources, per the commit's\n# citation-check pass.\n\n$d = $word.ActiveDocument\n\n$headingTitles = @(\"Introduction\", \"Overview of Microfinance\", \"Impact of Microfinance Institutions\", \"Effectiveness in Sustainable Development\", \"Conclusion\")\n\n$replacements = @(\n    [PSCustomObject]@{ Marker = \"Scholarly work on microfinance reveals a complex l\"; After = \"Scholarly work on microfinance reveals a complex landscape of both accomplishments and challenges. Success stories, such as that of the Grameen Bank, underscore the potential of microfinance institutions to empower marginalized populations by improving access to essential services and fostering economic activity within local communities (Smith). However, critiques highlight that while microfinance institutions have achieved economic success in operating microcredit, they often fall short of substantially raising the living standards of their clients (Smith). Furthermore, the literature points to persistent issues such as social divides and gender-based inequalities, which can hinder the full integration of women into economic activities (Smith). These critiques emphasize the need for addressing such challenges to ensure microfinance programs truly benefit the communities they aim to serve, thus promoting both financial and social inclusion (Smith).\" },\n    [PSCustomObject]@{ Marker = \"Microfinance institutions, exemplified by the Gram\"; After = \"Microfinance institutions, exemplified by the Grameen Bank, have played a transformative role in empowering local communities through financial inclusion and economic development. The Grameen Bank, renowned for its pioneering approach to microcredit, has significantly improved access to credit for impoverished individuals, thereby enabling them to invest in small-scale businesses and enhance their livelihoods (Ref-u503635). Another notable institution, such as BRAC, complements this impact by providing integrated services that include education and healthcare, further contributing to the well-being of community members (Ref-u503635). These institutions not only elevate household incomes but also stimulate broader economic activity by fostering entrepreneurship and creating employment opportunities in underserved regions (Ref-u503635). As a result, microfinance institutions have become pivotal in driving sustainable development and reducing poverty within local communities, despite facing persistent challenges that require ongoing attention (Ref-u503635).\" },\n    [PSCustomObject]@{ Marker = \"Furthermore, microfinance programs actively allevi\"; After = \"Furthermore, microfinance programs actively alleviate poverty through mechanisms like small loans and savings groups, which provide financial access to low-income individuals. These small loans empower recipients to start or expand businesses, leading to increased income and improved living standards (Ref-u026387). Savings groups, on the other hand, enable individuals to pool resources, creating a financial safety net that not only supports personal economic activities but also fosters a culture of saving and financial planning. The collective impact of these mechanisms enhances economic resilience among participants, allowing them to better withstand financial shocks and invest in opportunities that promote long-term stability (Ref-u026387). By integrating these approaches, microfinance programs play a critical role in poverty reduction, thereby contributing to broader economic development goals.\" },\n    [PSCustomObject]@{ Marker = \"However, despite their benefits, microfinance prog\"; After = \"However, despite their benefits, microfinance programs face several significant drawbacks that could impede their effectiveness. One major concern is the high interest rates charged by microfinance institutions, which are often necessary to cover operational costs but can burden borrowers and undermine the social mission of these programs (Ref-f532861). Additionally, the dependency created by microfinance services can lead to a cycle where borrowers become reliant on loans without achieving substantial improvements in their living standards (Ref-f532861). Gender-based challenges also persist, as women, who are often the primary beneficiaries of microfinance, may encounter social and cultural barriers that prevent them from fully participating in economic activities (Ref-f532861). These issues highlight the need for a more balanced approach that considers both the financial sustainability of microfinance institutions and the socio-economic empowerment of their clients, ensuring that the programs genuinely contribute to poverty alleviation and economic development.\" },\n    [PSCustomObject]@{ Marker = \"Microfinance has shown considerable potential in p\"; After = \"Microfinance has shown considerable potential in promoting sustainable development by improving the economic resilience of disadvantaged populations. Notably, microfinance institutions effectively alleviate credit constraints for impoverished communities, thereby enhancing financial inclusion and facilitating economic empowerment (Ref-u821407). These institutions provide critical financial resources that enable small businesses to flourish, which not only raises household incomes but also stimulates broader economic growth in underdeveloped regions (Ref-u821407). However, despite these positive outcomes, criticisms persist, particularly concerning the high interest rates and potential dependency issues that may limit the long-term impact of microfinance programs (Ref-u821407). Addressing these challenges requires the development of robust regulatory frameworks and policies that support the sustainability of microfinance initiatives while ensuring they contribute meaningfully to poverty reduction and sustainable development goals (Ref-u821407).\" }\n)\n\n$paragraphs = @($d.Paragraphs)\n\nforeach ($p in $paragraphs) {\n    $text = $p.Range.Text.Trim()\n    if ($headingTitles -contains $text -and $p.Style.NameLocal -eq \"Heading 2\") {\n        $p.Style = \"Normal\"\n    }\n}\n\nforeach ($p in $paragraphs) {\n    $text = $p.Range.Text\n    foreach ($r in $replacements) {\n        if ($text.StartsWith($r.Marker)) {\n            $p.Range.Text = $r.After\n            break\n        }\n    }\n}\n"}
